# Weekly update: insert a new week's data row for
# "Perejil" (Vega Modelo de Temuco) ahead of the existing series,
# pushing the prior rows (395-539) down by one row (396-540).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 395; existing rows 395-539 shift to 396-540.
$ws.Rows.Item(395).Insert()

# Populate the newly inserted row 395 with the new weekly observation.
$ws.Cells.Item(395, 1).Value  = 10
$ws.Cells.Item(395, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(395, 3).Value  = "La Araucanía"
$ws.Cells.Item(395, 4).Value  = 45215
$ws.Cells.Item(395, 5).Value  = 9
$ws.Cells.Item(395, 6).Value  = 100112044
$ws.Cells.Item(395, 7).Value  = "Perejil"
$ws.Cells.Item(395, 8).Value  = "Sin especificar"
$ws.Cells.Item(395, 9).Value  = "Primera"
$ws.Cells.Item(395, 10).Value = 55
$ws.Cells.Item(395, 11).Value = 5000
$ws.Cells.Item(395, 12).Value = 5000
$ws.Cells.Item(395, 13).Value = 5000
$ws.Cells.Item(395, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(395, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(395, 16).Value = 1667
$ws.Cells.Item(395, 17).Value = 3
$ws.Cells.Item(395, 18).Value = "Hortaliza"
